$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 533.6  # H80: 604.75 -> 533.6
$ws.Cells.Item(80, 9).Value = 526.5714  # I80: 573 -> 526.5714
$ws.Cells.Item(80, 10).Value = 550  # J80: 700 -> 550
$ws.Cells.Item(80, 11).Value = 1579.7142  # K80: 1719 -> 1579.7142
$ws.Cells.Item(80, 12).Value = 1650  # L80: 2100 -> 1650
$ws.Cells.Item(80, 13).Value = -581.7142000000001  # M80: -721 -> -581.7142000000001
$ws.Cells.Item(80, 14).Value = -3646  # N80: -4096 -> -3646

$ws.Cells.Item(83, 8).Value = 533.6  # H83: 604.75 -> 533.6
$ws.Cells.Item(83, 9).Value = 526.5714  # I83: 573 -> 526.5714
$ws.Cells.Item(83, 10).Value = 550  # J83: 700 -> 550
$ws.Cells.Item(83, 11).Value = 4739.1426  # K83: 5157 -> 4739.1426
$ws.Cells.Item(83, 12).Value = 4950  # L83: 6300 -> 4950
$ws.Cells.Item(83, 13).Value = 252.8573999999999  # M83: -165 -> 252.8573999999999
$ws.Cells.Item(83, 14).Value = -14934  # N83: -16284 -> -14934

$ws.Cells.Item(101, 8).Value = 403.5  # H101: 332 -> 403.5
$ws.Cells.Item(101, 9).Value = 338  # I101: 332 -> 338
$ws.Cells.Item(101, 10).Value = 600  # J101: 0 -> 600
$ws.Cells.Item(101, 11).Value = 1014  # K101: 996 -> 1014
$ws.Cells.Item(101, 12).Value = 1800  # L101: 0 -> 1800
$ws.Cells.Item(101, 13).Value = 608  # M101: 626 -> 608
$ws.Cells.Item(101, 14).Value = -5044  # N101: None -> -5044

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(28, 8).Value = 18470.5  # H28: 13470.667 -> 18470.5
$ws.Cells.Item(28, 9).Value = 18470.5  # I28: 13470.667 -> 18470.5
$ws.Cells.Item(28, 11).Value = 18470.5  # K28: 13470.667 -> 18470.5
$ws.Cells.Item(28, 13).Value = -18278.5  # M28: -13278.667 -> -18278.5

$ws.Cells.Item(61, 8).Value = 7194.9  # H61: 7496.9 -> 7194.9
$ws.Cells.Item(61, 9).Value = 4737.25  # I61: 5492.25 -> 4737.25
$ws.Cells.Item(61, 11).Value = 4737.25  # K61: 5492.25 -> 4737.25
$ws.Cells.Item(61, 13).Value = -4525.25  # M61: -5280.25 -> -4525.25

$ws.Cells.Item(74, 8).Value = 10349.25  # H74: 10349.75 -> 10349.25
$ws.Cells.Item(74, 9).Value = 11065.667  # I74: 12799.5 -> 11065.667
$ws.Cells.Item(74, 10).Value = 8200  # J74: 7900 -> 8200
$ws.Cells.Item(74, 11).Value = 11065.667  # K74: 12799.5 -> 11065.667
$ws.Cells.Item(74, 12).Value = 8200  # L74: 7900 -> 8200
$ws.Cells.Item(74, 13).Value = -10191.667  # M74: -11925.5 -> -10191.667
$ws.Cells.Item(74, 14).Value = -9948  # N74: -9648 -> -9948

$ws.Cells.Item(77, 8).Value = 10349.25  # H77: 10349.75 -> 10349.25
$ws.Cells.Item(77, 9).Value = 11065.667  # I77: 12799.5 -> 11065.667
$ws.Cells.Item(77, 10).Value = 8200  # J77: 7900 -> 8200
$ws.Cells.Item(77, 11).Value = 55328.335  # K77: 63997.5 -> 55328.335
$ws.Cells.Item(77, 12).Value = 41000  # L77: 39500 -> 41000
$ws.Cells.Item(77, 13).Value = -50960.335  # M77: -59629.5 -> -50960.335
$ws.Cells.Item(77, 14).Value = -49736  # N77: -48236 -> -49736

$ws.Cells.Item(97, 8).Value = 539.1429000000001  # H97: 530.05884 -> 539.1429000000001
$ws.Cells.Item(97, 9).Value = 575.2727  # I97: 576.63635 -> 575.2727
$ws.Cells.Item(97, 10).Value = 406.66666  # J97: 444.66666 -> 406.66666
$ws.Cells.Item(97, 11).Value = 575.2727  # K97: 576.63635 -> 575.2727
$ws.Cells.Item(97, 12).Value = 406.66666  # L97: 444.66666 -> 406.66666
$ws.Cells.Item(97, 13).Value = -79.27269999999999  # M97: -80.63634999999999 -> -79.27269999999999
$ws.Cells.Item(97, 14).Value = -1398.66666  # N97: -1436.66666 -> -1398.66666

$ws.Cells.Item(99, 8).Value = 18470.5  # H99: 13470.667 -> 18470.5
$ws.Cells.Item(99, 9).Value = 18470.5  # I99: 13470.667 -> 18470.5
$ws.Cells.Item(99, 11).Value = 18470.5  # K99: 13470.667 -> 18470.5
$ws.Cells.Item(99, 13).Value = -15475.5  # M99: -10475.667 -> -15475.5

$ws.Cells.Item(102, 8).Value = 3028.2632  # H102: 3043.3157 -> 3028.2632
$ws.Cells.Item(102, 9).Value = 1452.7142  # I102: 1473.1428 -> 1452.7142
$ws.Cells.Item(102, 11).Value = 1452.7142  # K102: 1473.1428 -> 1452.7142
$ws.Cells.Item(102, 13).Value = 169.2858000000001  # M102: 148.8571999999999 -> 169.2858000000001

$ws.Cells.Item(110, 8).Value = 471.9091  # H110: 377.2143 -> 471.9091
$ws.Cells.Item(110, 9).Value = 535.7778  # I110: 409.33334 -> 535.7778
$ws.Cells.Item(110, 11).Value = 535.7778  # K110: 409.33334 -> 535.7778
$ws.Cells.Item(110, 13).Value = 1509.2222  # M110: 1635.66666 -> 1509.2222

$ws.Cells.Item(136, 8).Value = 7194.9  # H136: 7496.9 -> 7194.9
$ws.Cells.Item(136, 9).Value = 4737.25  # I136: 5492.25 -> 4737.25
$ws.Cells.Item(136, 11).Value = 14211.75  # K136: 16476.75 -> 14211.75
$ws.Cells.Item(136, 13).Value = -11661.75  # M136: -13926.75 -> -11661.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(43, 8).Value = 220000  # H43: 0 -> 220000
$ws.Cells.Item(43, 10).Value = 220000  # J43: 0 -> 220000
$ws.Cells.Item(43, 12).Value = 220000  # L43: 0 -> 220000
$ws.Cells.Item(43, 14).Value = -220362  # N43: None -> -220362

$ws.Cells.Item(105, 8).Value = 1650.8462  # H105: 1787.1666 -> 1650.8462
$ws.Cells.Item(105, 9).Value = 1622.2222  # I105: 1823.125 -> 1622.2222
$ws.Cells.Item(105, 11).Value = 1622.2222  # K105: 1823.125 -> 1622.2222
$ws.Cells.Item(105, 13).Value = 124.7778000000001  # M105: -76.125 -> 124.7778000000001

$ws.Cells.Item(135, 8).Value = 128358.336  # H135: 175513 -> 128358.336
$ws.Cells.Item(135, 10).Value = 128358.336  # J135: 175513 -> 128358.336
$ws.Cells.Item(135, 12).Value = 128358.336  # L135: 175513 -> 128358.336
$ws.Cells.Item(135, 14).Value = -138498.336  # N135: -185653 -> -138498.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1381.5714  # H16: 48766.855 -> 1381.5714
$ws.Cells.Item(16, 9).Value = 1442.4546  # I16: 59808.293 -> 1442.4546
$ws.Cells.Item(16, 10).Value = 1158.3334  # J16: 1840.75 -> 1158.3334
$ws.Cells.Item(16, 11).Value = 1442.4546  # K16: 59808.293 -> 1442.4546
$ws.Cells.Item(16, 12).Value = 1158.3334  # L16: 1840.75 -> 1158.3334
$ws.Cells.Item(16, 13).Value = -1155.4546  # M16: -59521.293 -> -1155.4546
$ws.Cells.Item(16, 14).Value = -1732.3334  # N16: -2414.75 -> -1732.3334

$ws.Cells.Item(31, 8).Value = 7178.2383  # H31: 8658.9 -> 7178.2383
$ws.Cells.Item(31, 9).Value = 3650.5557  # I31: 4538 -> 3650.5557
$ws.Cells.Item(31, 10).Value = 9824  # J31: 10032.533 -> 9824
$ws.Cells.Item(31, 11).Value = 3650.5557  # K31: 4538 -> 3650.5557
$ws.Cells.Item(31, 12).Value = 9824  # L31: 10032.533 -> 9824
$ws.Cells.Item(31, 13).Value = -3355.5557  # M31: -4243 -> -3355.5557
$ws.Cells.Item(31, 14).Value = -10414  # N31: -10622.533 -> -10414

$ws.Cells.Item(34, 8).Value = 7178.2383  # H34: 8658.9 -> 7178.2383
$ws.Cells.Item(34, 9).Value = 3650.5557  # I34: 4538 -> 3650.5557
$ws.Cells.Item(34, 10).Value = 9824  # J34: 10032.533 -> 9824
$ws.Cells.Item(34, 11).Value = 3650.5557  # K34: 4538 -> 3650.5557
$ws.Cells.Item(34, 12).Value = 9824  # L34: 10032.533 -> 9824
$ws.Cells.Item(34, 13).Value = -3448.5557  # M34: -4336 -> -3448.5557
$ws.Cells.Item(34, 14).Value = -10228  # N34: -10436.533 -> -10228

$ws.Cells.Item(99, 8).Value = 1832.4  # H99: 1902.4166 -> 1832.4
$ws.Cells.Item(99, 9).Value = 1720.6666  # I99: 1785.8 -> 1720.6666
$ws.Cells.Item(99, 10).Value = 2000  # J99: 1985.7142 -> 2000
$ws.Cells.Item(99, 11).Value = 1720.6666  # K99: 1785.8 -> 1720.6666
$ws.Cells.Item(99, 12).Value = 2000  # L99: 1985.7142 -> 2000
$ws.Cells.Item(99, 13).Value = -222.6666  # M99: -287.8 -> -222.6666
$ws.Cells.Item(99, 14).Value = -4996  # N99: -4981.7142 -> -4996

$ws.Cells.Item(107, 8).Value = 323  # H107: 309.4737 -> 323
$ws.Cells.Item(107, 9).Value = 218.66667  # I107: 206.92308 -> 218.66667
$ws.Cells.Item(107, 11).Value = 218.66667  # K107: 206.92308 -> 218.66667
$ws.Cells.Item(107, 13).Value = 1701.33333  # M107: 1713.07692 -> 1701.33333

$ws.Cells.Item(113, 8).Value = 1381.5714  # H113: 48766.855 -> 1381.5714
$ws.Cells.Item(113, 9).Value = 1442.4546  # I113: 59808.293 -> 1442.4546
$ws.Cells.Item(113, 10).Value = 1158.3334  # J113: 1840.75 -> 1158.3334
$ws.Cells.Item(113, 11).Value = 1442.4546  # K113: 59808.293 -> 1442.4546
$ws.Cells.Item(113, 12).Value = 1158.3334  # L113: 1840.75 -> 1158.3334
$ws.Cells.Item(113, 13).Value = 727.5454  # M113: -57638.293 -> 727.5454
$ws.Cells.Item(113, 14).Value = -5498.3334  # N113: -6180.75 -> -5498.3334

$ws.Cells.Item(126, 8).Value = 1832.4  # H126: 1902.4166 -> 1832.4
$ws.Cells.Item(126, 9).Value = 1720.6666  # I126: 1785.8 -> 1720.6666
$ws.Cells.Item(126, 10).Value = 2000  # J126: 1985.7142 -> 2000
$ws.Cells.Item(126, 11).Value = 5161.9998  # K126: 5357.4 -> 5161.9998
$ws.Cells.Item(126, 12).Value = 6000  # L126: 5957.142599999999 -> 6000
$ws.Cells.Item(126, 13).Value = -2691.9998  # M126: -2887.4 -> -2691.9998
$ws.Cells.Item(126, 14).Value = -10940  # N126: -10897.1426 -> -10940

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(52, 8).Value = 975  # H52: 1225 -> 975
$ws.Cells.Item(52, 10).Value = 975  # J52: 1225 -> 975
$ws.Cells.Item(52, 12).Value = 2925  # L52: 3675 -> 2925
$ws.Cells.Item(52, 14).Value = -3457  # N52: -4207 -> -3457

$ws.Cells.Item(56, 8).Value = 0  # H56: 4000 -> 0
$ws.Cells.Item(56, 9).Value = 0  # I56: 4000 -> 0
$ws.Cells.Item(56, 11).Value = 0  # K56: 4000 -> 0
$ws.Cells.Item(56, 13).Value = $null  # M56: -3470 -> (cleared)

$ws.Cells.Item(114, 8).Value = 485.7  # H114: 416.85715 -> 485.7
$ws.Cells.Item(114, 9).Value = 198  # I114: 153 -> 198
$ws.Cells.Item(114, 10).Value = 1157  # J114: 2000 -> 1157
$ws.Cells.Item(114, 11).Value = 594  # K114: 459 -> 594
$ws.Cells.Item(114, 12).Value = 3471  # L114: 6000 -> 3471
$ws.Cells.Item(114, 13).Value = 2660  # M114: 2795 -> 2660
$ws.Cells.Item(114, 14).Value = -9979  # N114: -12508 -> -9979

$ws.Cells.Item(117, 8).Value = 590.3333  # H117: 601 -> 590.3333
$ws.Cells.Item(117, 9).Value = 74.5  # I117: 122.5 -> 74.5
$ws.Cells.Item(117, 11).Value = 223.5  # K117: 367.5 -> 223.5
$ws.Cells.Item(117, 13).Value = 3218.5  # M117: 3074.5 -> 3218.5

$ws.Cells.Item(121, 8).Value = 1749.75  # H121: 4000 -> 1749.75
$ws.Cells.Item(121, 9).Value = 999.5  # I121: 0 -> 999.5
$ws.Cells.Item(121, 10).Value = 2500  # J121: 4000 -> 2500
$ws.Cells.Item(121, 11).Value = 2998.5  # K121: 0 -> 2998.5
$ws.Cells.Item(121, 12).Value = 7500  # L121: 12000 -> 7500
$ws.Cells.Item(121, 13).Value = -1688.5  # M121: None -> -1688.5
$ws.Cells.Item(121, 14).Value = -10120  # N121: -14620 -> -10120

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 327.3889  # H97: 450.6842 -> 327.3889
$ws.Cells.Item(97, 9).Value = 292.93332  # I97: 301.07144 -> 292.93332
$ws.Cells.Item(97, 10).Value = 499.66666  # J97: 869.6 -> 499.66666
$ws.Cells.Item(97, 11).Value = 292.93332  # K97: 301.07144 -> 292.93332
$ws.Cells.Item(97, 12).Value = 499.66666  # L97: 869.6 -> 499.66666
$ws.Cells.Item(97, 13).Value = 203.06668  # M97: 194.92856 -> 203.06668
$ws.Cells.Item(97, 14).Value = -1491.66666  # N97: -1861.6 -> -1491.66666

$ws.Cells.Item(101, 8).Value = 45995  # H101: 23665 -> 45995
$ws.Cells.Item(101, 10).Value = 45995  # J101: 23665 -> 45995
$ws.Cells.Item(101, 12).Value = 45995  # L101: 23665 -> 45995
$ws.Cells.Item(101, 14).Value = -52485  # N101: -30155 -> -52485

$ws.Cells.Item(113, 8).Value = 6254.2144  # H113: 6092.6 -> 6254.2144
$ws.Cells.Item(113, 9).Value = 4687.3335  # I113: 4564.857 -> 4687.3335
$ws.Cells.Item(113, 11).Value = 4687.3335  # K113: 4564.857 -> 4687.3335
$ws.Cells.Item(113, 13).Value = -2517.3335  # M113: -2394.857 -> -2517.3335

$ws.Cells.Item(126, 8).Value = 2866.6667  # H126: 5000 -> 2866.6667
$ws.Cells.Item(126, 9).Value = 100  # I126: 5000 -> 100
$ws.Cells.Item(126, 10).Value = 4250  # J126: 0 -> 4250
$ws.Cells.Item(126, 11).Value = 300  # K126: 15000 -> 300
$ws.Cells.Item(126, 12).Value = 12750  # L126: 0 -> 12750
$ws.Cells.Item(126, 13).Value = 2170  # M126: -12530 -> 2170
$ws.Cells.Item(126, 14).Value = -17690  # N126: None -> -17690

$ws.Cells.Item(132, 8).Value = 2229  # H132: 1182.75 -> 2229
$ws.Cells.Item(132, 9).Value = 2229  # I132: 1182.3334 -> 2229
$ws.Cells.Item(132, 10).Value = 0  # J132: 1184 -> 0
$ws.Cells.Item(132, 11).Value = 6687  # K132: 3547.0002 -> 6687
$ws.Cells.Item(132, 12).Value = 0  # L132: 3552 -> 0
$ws.Cells.Item(132, 13).Value = $null  # M132: -1017.0002 -> (cleared)
$ws.Cells.Item(132, 14).Value = $null  # N132: -8612 -> (cleared)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 5000  # H2: 4332.6665 -> 5000
$ws.Cells.Item(2, 9).Value = 5000  # I2: 4499 -> 5000
$ws.Cells.Item(2, 10).Value = 0  # J2: 4000 -> 0
$ws.Cells.Item(2, 11).Value = 5000  # K2: 4499 -> 5000
$ws.Cells.Item(2, 12).Value = 0  # L2: 4000 -> 0
$ws.Cells.Item(2, 13).Value = $null  # M2: -4387 -> (cleared)
$ws.Cells.Item(2, 14).Value = $null  # N2: -4224 -> (cleared)

$ws.Cells.Item(100, 8).Value = 5968.385  # H100: 6694.8335 -> 5968.385
$ws.Cells.Item(100, 9).Value = 1931.5  # I100: 2667.8 -> 1931.5
$ws.Cells.Item(100, 10).Value = 9428.571  # J100: 9571.286 -> 9428.571
$ws.Cells.Item(100, 11).Value = 1931.5  # K100: 2667.8 -> 1931.5
$ws.Cells.Item(100, 12).Value = 9428.571  # L100: 9571.286 -> 9428.571
$ws.Cells.Item(100, 13).Value = -1390.5  # M100: -2126.8 -> -1390.5
$ws.Cells.Item(100, 14).Value = -10510.571  # N100: -10653.286 -> -10510.571

$ws.Cells.Item(136, 8).Value = 3002.5625  # H136: 2984.625 -> 3002.5625
$ws.Cells.Item(136, 9).Value = 2788  # I136: 2850.3333 -> 2788
$ws.Cells.Item(136, 10).Value = 3932.3333  # J136: 4999 -> 3932.3333
$ws.Cells.Item(136, 11).Value = 8364  # K136: 8550.999899999999 -> 8364
$ws.Cells.Item(136, 12).Value = 11796.9999  # L136: 14997 -> 11796.9999
$ws.Cells.Item(136, 13).Value = -5814  # M136: -6000.999899999999 -> -5814
$ws.Cells.Item(136, 14).Value = -16896.9999  # N136: -20097 -> -16896.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1086.0667  # H100: 1135.6428 -> 1086.0667
$ws.Cells.Item(100, 9).Value = 899.25  # I100: 945.36365 -> 899.25
$ws.Cells.Item(100, 11).Value = 1798.5  # K100: 1890.7273 -> 1798.5
$ws.Cells.Item(100, 13).Value = -1257.5  # M100: -1349.7273 -> -1257.5

$ws.Cells.Item(107, 8).Value = 392.13333  # H107: 414.7143 -> 392.13333
$ws.Cells.Item(107, 9).Value = 392.13333  # I107: 414.7143 -> 392.13333
$ws.Cells.Item(107, 11).Value = 1176.39999  # K107: 1244.1429 -> 1176.39999
$ws.Cells.Item(107, 13).Value = 743.6000100000001  # M107: 675.8571000000002 -> 743.6000100000001
